$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -8.056686313420776
$ws.Range("J25").Value = -7.844789947354293
$ws.Range("K25").Value = -2.492288476983724
$ws.Range("I26").Value = -7.895668583223774
$ws.Range("J26").Value = -2.543167112853204
$ws.Range("K26").Value = -2.64959566916621
$ws.Range("H27").Value = -7.848163674882571
$ws.Range("I27").Value = -2.495662204512
$ws.Range("J27").Value = -2.602090760825007
$ws.Range("K27").Value = -0.8533747632109219
$ws.Range("G28").Value = -7.838666639181375
$ws.Range("H28").Value = -2.486165168810806
$ws.Range("I28").Value = -2.592593725123812
$ws.Range("J28").Value = -0.8438777275097269
$ws.Range("K28").Value = 1.781949470911521
$ws.Range("F29").Value = -8.107842639081671
$ws.Range("G29").Value = -2.755341168711102
$ws.Range("H29").Value = -2.861769725024108
$ws.Range("I29").Value = -1.113053727410023
$ws.Range("J29").Value = 1.512773471011225
$ws.Range("K29").Value = -4.130805822877837
$ws.Range("E30").Value = -8.155936619030644
$ws.Range("F30").Value = -2.803435148660075
$ws.Range("G30").Value = -2.909863704973081
$ws.Range("H30").Value = -1.161147707358996
$ws.Range("I30").Value = 1.464679491062252
$ws.Range("J30").Value = -4.17889980282681
$ws.Range("K30").Value = -0.02515630166823679
$ws.Range("D31").Value = -8.355093825182436
$ws.Range("E31").Value = -3.002592354811866
$ws.Range("F31").Value = -3.109020911124873
$ws.Range("G31").Value = -1.360304913510788
$ws.Range("H31").Value = 1.26552228491046
$ws.Range("I31").Value = -4.378057008978601
$ws.Range("J31").Value = -0.2243135078200282
$ws.Range("K31").Value = -0.8901600277517188
$ws.Range("C32").Value = -8.642971246787017
$ws.Range("D32").Value = -3.290469776416447
$ws.Range("E32").Value = -3.396898332729454
$ws.Range("F32").Value = -1.648182335115369
$ws.Range("G32").Value = 0.9776448633058792
$ws.Range("H32").Value = -4.665934430583182
$ws.Range("I32").Value = -0.5121909294246094
$ws.Range("J32").Value = -1.1780374493563
$ws.Range("K32").Value = -1.935101045603123
$ws.Range("B33").Value = -10.45854990725616
$ws.Range("C33").Value = -5.106048436885593
$ws.Range("D33").Value = -5.212476993198599
$ws.Range("E33").Value = -3.463760995584514
$ws.Range("F33").Value = -0.8379337971632659
$ws.Range("G33").Value = -6.481513091052327
$ws.Range("H33").Value = -2.327769589893754
$ws.Range("I33").Value = -2.993616109825445
$ws.Range("J33").Value = -3.750679706072268
$ws.Range("K33").Value = -0.06145210016119695
$ws.Range("B34").Value = -4.146418257367844
$ws.Range("C34").Value = -4.252846813680851
$ws.Range("D34").Value = -2.504130816066766
$ws.Range("E34").Value = 0.1216963823544824
$ws.Range("F34").Value = -5.521882911534579
$ws.Range("G34").Value = -1.368139410376006
$ws.Range("H34").Value = -2.033985930307697
$ws.Range("I34").Value = -2.79104952655452
$ws.Range("J34").Value = 0.8981780793565513
$ws.Range("K34").Value = 0.4649817484139299
$ws.Range("B35").Value = -4.065174135534914
$ws.Range("C35").Value = -2.316458137920828
$ws.Range("D35").Value = 0.3093690605004195
$ws.Range("E35").Value = -5.334210233388642
$ws.Range("F35").Value = -1.180466732230069
$ws.Range("G35").Value = -1.84631325216176
$ws.Range("H35").Value = -2.603376848408582
$ws.Range("I35").Value = 1.085850757502488
$ws.Range("J35").Value = 0.6526544265598669
$ws.Range("K35").Value = -0.2441609490878097
$ws.Range("B36").Value = -2.079249613080133
$ws.Range("C36").Value = 0.5465775853411154
$ws.Range("D36").Value = -5.097001708547946
$ws.Range("E36").Value = -0.9432582073893732
$ws.Range("F36").Value = -1.609104727321064
$ws.Range("G36").Value = -2.366168323567887
$ws.Range("H36").Value = 1.323059282343184
$ws.Range("I36").Value = 0.8898629514005629
$ws.Range("J36").Value = -0.006952424247113734
$ws.Range("K36").Value = 1.259738135287552
$ws.Range("B37").Value = 0.7313210777328913
$ws.Range("C37").Value = -4.91225821615617
$ws.Range("D37").Value = -0.7585147149975973
$ws.Range("E37").Value = -1.424361234929288
$ws.Range("F37").Value = -2.18142483117611
$ws.Range("G37").Value = 1.50780277473496
$ws.Range("H37").Value = 1.074606443792339
$ws.Range("I37").Value = 0.1777910681446622
$ws.Range("J37").Value = 1.444481627679328
$ws.Range("K37").Value = 1.851502671729013
$ws.Range("B38").Value = -5.105156740230972
$ws.Range("C38").Value = -0.9514132390723986
$ws.Range("D38").Value = -1.617259759004089
$ws.Range("E38").Value = -2.374323355250912
$ws.Range("F38").Value = 1.314904250660159
$ws.Range("G38").Value = 0.8817079197175375
$ws.Range("H38").Value = -0.01510745593013913
$ws.Range("I38").Value = 1.251583103604527
$ws.Range("J38").Value = 1.658604147654212
$ws.Range("K38").Value = -0.6382024769694439
$ws.Range("B39").Value = -0.3421915800502259
$ws.Range("C39").Value = -1.008038099981917
$ws.Range("D39").Value = -1.765101696228739
$ws.Range("E39").Value = 1.924125909682332
$ws.Range("F39").Value = 1.49092957873971
$ws.Range("G39").Value = 0.5941142030920336
$ws.Range("H39").Value = 1.860804762626699
$ws.Range("I39").Value = 2.267825806676385
$ws.Range("J39").Value = -0.02898081794727114
$ws.Range("K39").Value = 0.5598170418495986
$ws.Range("B40").Value = -0.9817114683876786
$ws.Range("C40").Value = -1.738775064634501
$ws.Range("D40").Value = 1.950452541276569
$ws.Range("E40").Value = 1.517256210333948
$ws.Range("F40").Value = 0.6204408346862713
$ws.Range("G40").Value = 1.887131394220937
$ws.Range("H40").Value = 2.294152438270622
$ws.Range("I40").Value = -0.002654186353033383
$ws.Range("J40").Value = 0.5861436734438363
$ws.Range("K40").Value = 0.4572641887935098
$ws.Range("B41").Value = -1.640936427035204
$ws.Range("C41").Value = 2.048291178875867
$ws.Range("D41").Value = 1.615094847933246
$ws.Range("E41").Value = 0.718279472285569
$ws.Range("F41").Value = 1.984970031820235
$ws.Range("G41").Value = 2.39199107586992
$ws.Range("H41").Value = 0.09518445124626429
$ws.Range("I41").Value = 0.683982311043134
$ws.Range("J41").Value = 0.5551028263928075
$ws.Range("K41").Value = 0.1623546509229741
$ws.Range("B42").Value = 2.415062601821465
$ws.Range("C42").Value = 1.981866270878843
$ws.Range("D42").Value = 1.085050895231167
$ws.Range("E42").Value = 2.351741454765833
$ws.Range("F42").Value = 2.758762498815518
$ws.Range("G42").Value = 0.4619558741918621
$ws.Range("H42").Value = 1.050753733988732
$ws.Range("I42").Value = 0.9218742493384053
$ws.Range("J42").Value = 0.5291260738685719
$ws.Range("K42").Value = 0.002169620121424964
$ws.Range("B43").Value = 3.478476033585537
$ws.Range("C43").Value = 2.58166065793786
$ws.Range("D43").Value = 3.848351217472526
$ws.Range("E43").Value = 4.255372261522211
$ws.Range("F43").Value = 1.958565636898555
$ws.Range("G43").Value = 2.547363496695425
$ws.Range("H43").Value = 2.418484012045099
$ws.Range("I43").Value = 2.025735836575265
$ws.Range("J43").Value = 1.498779382828118
$ws.Range("K43").Value = 0.9313234029264381
$ws.Range("B44").Value = 1.584099193749569
$ws.Range("C44").Value = 2.850789753284235
$ws.Range("D44").Value = 3.25781079733392
$ws.Range("E44").Value = 0.9610041727102643
$ws.Range("F44").Value = 1.549802032507134
$ws.Range("G44").Value = 1.420922547856807
$ws.Range("H44").Value = 1.028174372386974
$ws.Range("I44").Value = 0.5012179186398271
$ws.Range("J44").Value = -0.06623806126185294
$ws.Range("B45").Value = 2.466562918591728
$ws.Range("C45").Value = 2.873583962641413
$ws.Range("D45").Value = 0.5767773380177572
$ws.Range("E45").Value = 1.165575197814627
$ws.Range("F45").Value = 1.0366957131643
$ws.Range("G45").Value = 0.643947537694467
$ws.Range("H45").Value = 0.1169910839473201
$ws.Range("I45").Value = -0.45046489595436
$ws.Range("B46").Value = 2.246662553384068
$ws.Range("C46").Value = -0.05014407123958764
$ws.Range("D46").Value = 0.5386537885572821
$ws.Range("E46").Value = 0.4097743039069556
$ws.Range("F46").Value = 0.01702612843712215
$ws.Range("G46").Value = -0.5099303253100248
$ws.Range("H46").Value = -1.077386305211705
$ws.Range("B47").Value = -0.4928110405894273
$ws.Range("C47").Value = 0.09598681920744241
$ws.Range("D47").Value = -0.03289266544288409
$ws.Range("E47").Value = -0.4256408409127175
$ws.Range("F47").Value = -0.9525972946598644
$ws.Range("G47").Value = -1.520053274561544
$ws.Range("B48").Value = 0.2881944199355046
$ws.Range("C48").Value = 0.1593149352851781
$ws.Range("D48").Value = -0.2334332401846553
$ws.Range("E48").Value = -0.7603896939318022
$ws.Range("F48").Value = -1.327845673833482
$ws.Range("B49").Value = -0.2265482812103926
$ws.Range("C49").Value = -0.6192964566802259
$ws.Range("D49").Value = -1.146252910427373
$ws.Range("E49").Value = -1.713708890329053
$ws.Range("B50").Value = -0.5033433667621097
$ws.Range("C50").Value = -1.030299820509257
$ws.Range("D50").Value = -1.597755800410937
$ws.Range("B51").Value = -1.152650332186175
$ws.Range("C51").Value = -1.720106312087855
$ws.Range("B52").Value = -2.040355454506476
